# StructureDefinition-match-period.xlsx update
# - bump Version 5.0.0 -> 6.0.0
# - bump Date to the new publication timestamp
# - replace the broken "Contact / No display for ContactDetail" row with a
#   proper Publisher value, and add a new Jurisdiction row
# - remove the old duplicated "Contact" row (sheet shrinks from 21 to 20 rows)
# - update the root Extension's Short/Definition text on the Elements sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value  = "6.0.0"
$ws1.Range("B8").Value  = "2022-01-21T20:46:54+00:00"
$ws1.Range("B9").Value  = "Alvearie Team"

$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 is still the stale duplicate "Contact" / "No display for ContactDetail"
# row at this point; delete it so everything below shifts up one row.
$ws1.Rows.Item(11).Delete()

$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Match Period"
$ws2.Range("L2").Value = "Time period when match is valid"
